# Alvearie FHIR IG deploy update:
#  - StructureDefinition version bump 5.0.0 -> 6.0.0
#  - Date refresh
#  - Publisher/Jurisdiction metadata filled in (replacing the duplicated
#    "Contact" row with a real Publisher value + a new Jurisdiction row)
#  - Root "Extension" element's Short/Definition updated to reflect the
#    specific extension (Claim Status) instead of the generic placeholder

$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" --------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version
$meta.Range("B3").Value = "6.0.0"

# Date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a real value
$meta.Range("B9").Value = "Alvearie Team"

# Row 11 ("Contact" / "No display for ContactDetail") was a duplicate of
# row 10 - remove it so the rows below shift up by one.
$meta.Rows.Item(11).Delete()

# The (now single) row 10 becomes the new "Jurisdiction" property.
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Sheet "Elements" --------------------------------------------------
$elem = $wb.Worksheets.Item("Elements")

# Root Extension row: Short / Definition now reflect this specific
# extension rather than the generic boilerplate text.
$elem.Range("K2").Value = "Claim Status"
$elem.Range("L2").Value = "Code for the status of a claim"
